# Automatische test-sync: 2025-06-19 17:55:30
# Appends the two new "Klacht over levering" mail-log rows to the Logs
# sheet, extends the conditional formatting ranges to cover them, and
# refreshes the "Klacht" count on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# --- Append the two new log rows (28 and 29) -----------------------------

$wsLogs.Range("A28").Value = "Klacht over levering"
$wsLogs.Range("B28").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C28").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$wsLogs.Range("D28").Value = "Klacht"
$wsLogs.Range("F28").Value = "2025-06-19 17:55:20"
$wsLogs.Range("G28").Value = "Nee"

$wsLogs.Range("A29").Value = "Klacht over levering"
$wsLogs.Range("B29").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C29").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$wsLogs.Range("D29").Value = "Klacht"
$wsLogs.Range("F29").Value = "2025-06-19 17:55:21"
$wsLogs.Range("G29").Value = "Nee"

# --- Extend the conditional formatting sqref ranges to include rows 28-29 -

$fcCategorie = $wsLogs.Range("D2:D27").FormatConditions.Item(1)
$fcCategorie.ModifyAppliesToRange($wsLogs.Range("D2:D29"))

$fcBeantwoord = $wsLogs.Range("G2:G27").FormatConditions.Item(1)
$fcBeantwoord.ModifyAppliesToRange($wsLogs.Range("G2:G29"))

# --- Refresh the Dashboard "Klacht" tally (3 -> 5) ------------------------

$wsDash.Range("B4").Value = 5
